$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the rule action text in G7: the "Modified Rating Factor" rule now
#    multiplies the factor and rating factor together instead of taking $param.
$ws.Range("G7").Value = "`$p.setModifiedRatingFactor(`$p.getFactor()*`$p.getRatingFactor());"

# 2. The "Modified Rating Factor" column (G9:G18) no longer holds the
#    pre-computed E*F formula/value - clear it out while keeping the cell
#    formatting/style intact.
$ws.Range("G9:G18").ClearContents()

# 3. The author's cursor/selection ended up on F26 when the file was saved.
$ws.Range("F26").Select()
